# Refresh the cryptocurrency price/volume snapshot (GitHub Actions daily update).
# Note: several "Price" values are numeric-looking (e.g. "555.95", "0.999", "1.00")
# and a couple of rows (38/39 and 43-46) change which coin/link occupies a given
# row (ranking reshuffle), so those are rewritten explicitly too. A leading
# apostrophe forces Excel to keep a numeric-looking price string as text,
# matching the workbook's existing inline-string formatting for that column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.550.64'
$ws.Range("E2").Value = '  -5.62%  '

$ws.Range("D3").Value = '3.303.62'
$ws.Range("E3").Value = '  -6.66%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = '''555.95'
$ws.Range("E5").Value = '  -4.96%  '

$ws.Range("D6").Value = '''179.54'
$ws.Range("E6").Value = '  -7.99%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '''0.588'
$ws.Range("E8").Value = '  -3.57%  '

$ws.Range("D9").Value = '3.296.53'
$ws.Range("E9").Value = '  -6.54%  '

$ws.Range("D10").Value = '''0.184'
$ws.Range("E10").Value = '  -10.65%  '

$ws.Range("D11").Value = '''0.585'
$ws.Range("E11").Value = '  -6.88%  '

$ws.Range("D12").Value = '''47.03'
$ws.Range("E12").Value = '  -10.23%  '

$ws.Range("D13").Value = '''0.0000263'
$ws.Range("E13").Value = '  -8.36%  '

$ws.Range("D14").Value = '3.842.36'
$ws.Range("E14").Value = '  -6.29%  '

$ws.Range("D15").Value = '''8.50'
$ws.Range("E15").Value = '  -7.47%  '

$ws.Range("D16").Value = '''597.61'
$ws.Range("E16").Value = '  -9.85%  '

$ws.Range("D17").Value = '''18.00'
$ws.Range("E17").Value = '  -2.03%  '

$ws.Range("D18").Value = '65.599.61'
$ws.Range("E18").Value = '  -5.65%  '

$ws.Range("E19").Value = '  -4.00%  '

$ws.Range("D20").Value = '3.302.72'
$ws.Range("E20").Value = '  -6.52%  '

$ws.Range("D21").Value = '''11.37'
$ws.Range("E21").Value = '  -8.74%  '

$ws.Range("D22").Value = '''0.898'
$ws.Range("E22").Value = '  -6.50%  '

$ws.Range("D23").Value = '''17.29'
$ws.Range("E23").Value = '  -4.22%  '

$ws.Range("D24").Value = '''102.57'
$ws.Range("E24").Value = '  -1.40%  '

$ws.Range("D25").Value = '''5.02'
$ws.Range("E25").Value = '  -6.91%  '

$ws.Range("D26").Value = '''3.95'

$ws.Range("E27").Value = '  -0.50%  '

$ws.Range("D28").Value = '''2.66'
$ws.Range("E28").Value = '  -8.90%  '

$ws.Range("D29").Value = '''9.26'
$ws.Range("E29").Value = '  -8.18%  '

$ws.Range("D30").Value = '''8.64'
$ws.Range("E30").Value = '  -9.39%  '

$ws.Range("D31").Value = '''30.40'
$ws.Range("E31").Value = '  -8.32%  '

$ws.Range("D32").Value = '''3.86'
$ws.Range("E32").Value = '  -10.95%  '

$ws.Range("D33").Value = '''6.21'
$ws.Range("E33").Value = '  -7.80%  '

$ws.Range("D34").Value = '''10.99'
$ws.Range("E34").Value = '  -6.33%  '

$ws.Range("E35").Value = '  -6.13%  '

$ws.Range("D36").Value = '3.793.83'
$ws.Range("E36").Value = '  +0.79%  '

$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '''519.43'
$ws.Range("E38").Value = '  +4.65%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '''56.07'
$ws.Range("E39").Value = '  -9.01%  '

$ws.Range("D40").Value = '''3.46'
$ws.Range("E40").Value = '  -8.83%  '

$ws.Range("E41").Value = '  -12.63%  '

$ws.Range("D42").Value = '''2.62'
$ws.Range("E42").Value = '  -8.67%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").Value = '''0.123'
$ws.Range("E43").Value = '  -7.92%  '

$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''31.77'
$ws.Range("E44").Value = '  -8.05%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '''0.337'
$ws.Range("E45").Value = '  -8.87%  '

$ws.Range("B46").Value = 'CoreDAO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D46").Value = '''3.19'
$ws.Range("E46").Value = '  +18.09%  '

$ws.Range("D47").Value = '''3.27'
$ws.Range("E47").Value = '  -3.93%  '

$ws.Range("E48").Value = '  -8.46%  '

$ws.Range("E49").Value = '  -4.97%  '

$ws.Range("E50").Value = '  -9.71%  '

$ws.Range("D51").Value = '''1.00'
$ws.Range("E51").Value = '  -0.11%  '
